# Scheduled runner update: refresh market-board price/profit columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) for the affected leve rows across the ALC, ARM, CRP,
# CUL, GSM and LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 173.83333
$ws.Range("I5").Value = 147
$ws.Range("J5").Value = 200.66667
$ws.Range("K5").Value = 147
$ws.Range("L5").Value = 200.66667
$ws.Range("M5").Value = -32
$ws.Range("N5").Value = -430.66667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 738489.75
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 738489.75
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2215469.25
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2215805.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2005.909
$ws.Range("I40").Value = 1908.8889
$ws.Range("J40").Value = 2073.077
$ws.Range("K40").Value = 1908.8889
$ws.Range("L40").Value = 2073.077
$ws.Range("M40").Value = -1733.8889
$ws.Range("N40").Value = -2423.077

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1993.2916
$ws.Range("I70").Value = 2618.182
$ws.Range("J70").Value = 1464.5385
$ws.Range("K70").Value = 7854.545999999999
$ws.Range("L70").Value = 4393.6155
$ws.Range("M70").Value = -7584.545999999999
$ws.Range("N70").Value = -4933.6155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1993.2916
$ws.Range("I73").Value = 2618.182
$ws.Range("J73").Value = 1464.5385
$ws.Range("K73").Value = 7854.545999999999
$ws.Range("L73").Value = 4393.6155
$ws.Range("M73").Value = -6918.545999999999
$ws.Range("N73").Value = -6265.6155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1451.8966
$ws.Range("I86").Value = 1741.4
$ws.Range("J86").Value = 1141.7142
$ws.Range("K86").Value = 1741.4
$ws.Range("L86").Value = 1141.7142
$ws.Range("M86").Value = -618.4000000000001
$ws.Range("N86").Value = -3387.7142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1451.8966
$ws.Range("I89").Value = 1741.4
$ws.Range("J89").Value = 1141.7142
$ws.Range("K89").Value = 8707
$ws.Range("L89").Value = 5708.571
$ws.Range("M89").Value = -3091
$ws.Range("N89").Value = -16940.571

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 101264
$ws.Range("I100").Value = 126205
$ws.Range("J100").Value = 1500
$ws.Range("K100").Value = 126205
$ws.Range("L100").Value = 1500
$ws.Range("M100").Value = -125664
$ws.Range("N100").Value = -2582

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5541.9116
$ws.Range("I137").Value = 3636.3235
$ws.Range("J137").Value = 7447.5
$ws.Range("K137").Value = 10908.9705
$ws.Range("L137").Value = 22342.5
$ws.Range("M137").Value = -8358.970499999999
$ws.Range("N137").Value = -27442.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2007.1464
$ws.Range("I32").Value = 1824.7222
$ws.Range("J32").Value = 3320.6
$ws.Range("K32").Value = 1824.7222
$ws.Range("L32").Value = 3320.6
$ws.Range("M32").Value = -1537.7222
$ws.Range("N32").Value = -3894.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1553.4103
$ws.Range("I110").Value = 1249.4242
$ws.Range("J110").Value = 3225.3333
$ws.Range("K110").Value = 1249.4242
$ws.Range("L110").Value = 3225.3333
$ws.Range("M110").Value = 795.5758000000001
$ws.Range("N110").Value = -7315.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66.75
$ws.Range("I7").Value = 53.8
$ws.Range("J7").Value = 76
$ws.Range("K7").Value = 53.8
$ws.Range("L7").Value = 76
$ws.Range("M7").Value = 59.2
$ws.Range("N7").Value = -302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 623
$ws.Range("I22").Value = 681.625
$ws.Range("J22").Value = 466.66666
$ws.Range("K22").Value = 681.625
$ws.Range("L22").Value = 466.66666
$ws.Range("M22").Value = -331.625
$ws.Range("N22").Value = -1166.66666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5150.8335
$ws.Range("I31").Value = 4594.2856
$ws.Range("J31").Value = 5505
$ws.Range("K31").Value = 4594.2856
$ws.Range("L31").Value = 5505
$ws.Range("M31").Value = -4299.2856
$ws.Range("N31").Value = -6095

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5150.8335
$ws.Range("I34").Value = 4594.2856
$ws.Range("J34").Value = 5505
$ws.Range("K34").Value = 4594.2856
$ws.Range("L34").Value = 5505
$ws.Range("M34").Value = -4392.2856
$ws.Range("N34").Value = -5909

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 667333.7
$ws.Range("I22").Value = 1000500.5
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3001501.5
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = -3001332.5
$ws.Range("N22").Value = -3338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 667333.7
$ws.Range("I27").Value = 1000500.5
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 3001501.5
$ws.Range("L27").Value = 3000
$ws.Range("M27").Value = -3001399.5
$ws.Range("N27").Value = -3204

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 323500
$ws.Range("I32").Value = 414285.72
$ws.Range("J32").Value = 5750
$ws.Range("K32").Value = 1242857.16
$ws.Range("L32").Value = 17250
$ws.Range("M32").Value = -1242574.16
$ws.Range("N32").Value = -17816

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 912.95
$ws.Range("I68").Value = 675.3171
$ws.Range("J68").Value = 1078.0847
$ws.Range("K68").Value = 2025.9513
$ws.Range("L68").Value = 3234.2541
$ws.Range("M68").Value = -1214.9513
$ws.Range("N68").Value = -4856.2541

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 912.95
$ws.Range("I71").Value = 675.3171
$ws.Range("J71").Value = 1078.0847
$ws.Range("K71").Value = 6077.8539
$ws.Range("L71").Value = 9702.7623
$ws.Range("M71").Value = -2021.8539
$ws.Range("N71").Value = -17814.7623

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 31.90909
$ws.Range("I2").Value = 19.857143
$ws.Range("J2").Value = 53
$ws.Range("K2").Value = 19.857143
$ws.Range("L2").Value = 53
$ws.Range("M2").Value = 93.14285699999999
$ws.Range("N2").Value = -279

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4895.5654
$ws.Range("I70").Value = 4840.1333
$ws.Range("J70").Value = 4999.5
$ws.Range("K70").Value = 4840.1333
$ws.Range("L70").Value = 4999.5
$ws.Range("M70").Value = -4570.1333
$ws.Range("N70").Value = -5539.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4895.5654
$ws.Range("I73").Value = 4840.1333
$ws.Range("J73").Value = 4999.5
$ws.Range("K73").Value = 4840.1333
$ws.Range("L73").Value = 4999.5
$ws.Range("M73").Value = -3904.1333
$ws.Range("N73").Value = -6871.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 361
$ws.Range("I16").Value = 361
$ws.Range("K16").Value = 361
$ws.Range("M16").Value = -191

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1334.2084
$ws.Range("I46").Value = 1521.5
$ws.Range("J46").Value = 1072
$ws.Range("K46").Value = 1521.5
$ws.Range("L46").Value = 1072
$ws.Range("M46").Value = -1333.5
$ws.Range("N46").Value = -1448
